# Update "想去人数" (F column) figures across the "展览", "演出", and
# "全部类型" worksheets to match the refreshed data snapshot.
# "本地生活" only contains a header row, so no data changes are required there.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 6455
$ws.Range("F3").Value = 2602
$ws.Range("F5").Value = 1307
$ws.Range("F7").Value = 3186
$ws.Range("F8").Value = 379
$ws.Range("F9").Value = 174
$ws.Range("F11").Value = 8062
$ws.Range("F12").Value = 418
$ws.Range("F13").Value = 74
$ws.Range("F16").Value = 287
$ws.Range("F18").Value = 53
$ws.Range("F20").Value = 308
$ws.Range("F21").Value = 10034
$ws.Range("F23").Value = 279
$ws.Range("F30").Value = 80
$ws.Range("F31").Value = 124
$ws.Range("F33").Value = 2058
$ws.Range("F35").Value = 28
$ws.Range("F36").Value = 2099
$ws.Range("F37").Value = 4022
$ws.Range("F38").Value = 250
$ws.Range("F40").Value = 2008
$ws.Range("F41").Value = 1212
$ws.Range("F42").Value = 136
$ws.Range("F43").Value = 294
$ws.Range("F44").Value = 201
$ws.Range("F46").Value = 86
$ws.Range("F47").Value = 74
$ws.Range("F48").Value = 76
$ws.Range("F49").Value = 50

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 167
$ws.Range("F6").Value = 35
$ws.Range("F13").Value = 23
$ws.Range("F20").Value = 18

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6455
$ws.Range("F3").Value = 2602
$ws.Range("F4").Value = 167
$ws.Range("F6").Value = 1307
$ws.Range("F8").Value = 3186
$ws.Range("F9").Value = 379
$ws.Range("F11").Value = 174
$ws.Range("F13").Value = 8062
$ws.Range("F14").Value = 418
$ws.Range("F15").Value = 74
$ws.Range("F18").Value = 287
$ws.Range("F19").Value = 53
$ws.Range("F21").Value = 308
$ws.Range("F22").Value = 10034
$ws.Range("F23").Value = 279
$ws.Range("F28").Value = 23
$ws.Range("F30").Value = 80
$ws.Range("F31").Value = 124
$ws.Range("F33").Value = 2058
$ws.Range("F35").Value = 2099
$ws.Range("F36").Value = 4022
$ws.Range("F37").Value = 250
$ws.Range("F39").Value = 2009
$ws.Range("F41").Value = 1212
$ws.Range("F42").Value = 136
$ws.Range("F43").Value = 294
$ws.Range("F44").Value = 201
$ws.Range("F46").Value = 86
$ws.Range("F47").Value = 74
$ws.Range("F48").Value = 76
$ws.Range("F49").Value = 50
